# Weekly fruit/vegetable price update: a new weekly observation is inserted
# at row 11 (pushing the existing rows 11-100 down to 12-101), and the new
# row is populated with this week's Cilantro price data for
# "Terminal La Palmera de La Serena".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11; this shifts rows 11..100 down
# to 12..101 (and the sheet's used range / dimension grows to A1:R101).
$ws.Rows(11).Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44532
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112040
$ws.Range("G11").Value = "Cilantro"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 3200
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1900
$ws.Range("N11").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 1267
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = "Hortaliza"
